$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "last updated" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 10 de Agosto de 2020 a las 11:17"

# --- Update country data rows (values in columns B..H) ---
# Row 6: India
$ws.Range("B6").Value = 2217645
$ws.Range("C6").Value = 3508
$ws.Range("D6").Value = 1536259
$ws.Range("E6").Value = 636887
$ws.Range("G6").Value = 33
$ws.Range("H6").Value = 44499

# Row 7: Rusia
$ws.Range("B7").Value = 892654
$ws.Range("C7").Value = 5118
$ws.Range("D7").Value = 696681
$ws.Range("E7").Value = 180972
$ws.Range("G7").Value = 70
$ws.Range("H7").Value = 15001

# Row 18: Banglades
$ws.Range("B18").Value = 260507
$ws.Range("C18").Value = 2907
$ws.Range("D18").Value = 150437
$ws.Range("E18").Value = 106632
$ws.Range("G18").Value = 39
$ws.Range("H18").Value = 3438

# Row 26: Indonesia
$ws.Range("B26").Value = 127083
$ws.Range("C26").Value = 1687
$ws.Range("D26").Value = 82236
$ws.Range("E26").Value = 39082
$ws.Range("G26").Value = 42
$ws.Range("H26").Value = 5765

# Row 34: Israel
$ws.Range("B34").Value = 83540
$ws.Range("C34").Value = 538
$ws.Range("D34").Value = 58934
$ws.Range("E34").Value = 24000
$ws.Range("G34").Value = 6
$ws.Range("H34").Value = 606

# Row 47: Singapur
$ws.Range("B47").Value = 55292
$ws.Range("C47").Value = 188
$ws.Range("E47").Value = 6350

# Row 49: Polonia
$ws.Range("B49").Value = 52410
$ws.Range("C49").Value = 619
$ws.Range("D49").Value = 36877
$ws.Range("E49").Value = 13724
$ws.Range("G49").Value = 2
$ws.Range("H49").Value = 1809

# Row 55: Armenia
$ws.Range("B55").Value = 40433
$ws.Range("C55").Value = 23
$ws.Range("D55").Value = 32616
$ws.Range("E55").Value = 7021
$ws.Range("G55").Value = 5
$ws.Range("H55").Value = 796

# Row 71: Austria
$ws.Range("B71").Value = 22106
$ws.Range("C71").Value = 73
$ws.Range("D71").Value = 20010
$ws.Range("E71").Value = 1373
$ws.Range("G71").Value = 2
$ws.Range("H71").Value = 723

# Row 73: El Salvador
$ws.Range("B73").Value = 20872
$ws.Range("C73").Value = 449
$ws.Range("D73").Value = 9720
$ws.Range("E73").Value = 10589
$ws.Range("G73").Value = 14
$ws.Range("H73").Value = 563

# Row 74: Chequia
$ws.Range("B74").Value = 18355
$ws.Range("C74").Value = 2
$ws.Range("E74").Value = 5180

# Row 107: Hungria
$ws.Range("B107").Value = 4731
$ws.Range("C107").Value = 35
$ws.Range("D107").Value = 3525
$ws.Range("E107").Value = 601
$ws.Range("G107").Value = 3
$ws.Range("H107").Value = 605

# Row 111: Hong Kong
$ws.Range("B111").Value = 4149
$ws.Range("C111").Value = 69
$ws.Range("D111").Value = 2916
$ws.Range("E111").Value = 1178
$ws.Range("H111").Value = 55

# Row 122: Sri Lanka
$ws.Range("D122").Value = 2593
$ws.Range("E122").Value = 240

# Row 123: Eslovaquia
$ws.Range("B123").Value = 2599
$ws.Range("C123").Value = 3
$ws.Range("D123").Value = 1866
$ws.Range("E123").Value = 702

# Row 128: Lituania
$ws.Range("B128").Value = 2265
$ws.Range("C128").Value = 13
$ws.Range("E128").Value = 514

# Row 130: Estonia
$ws.Range("B130").Value = 2158
$ws.Range("C130").Value = 6
$ws.Range("D130").Value = 1962
$ws.Range("E130").Value = 133

# --- Row 148/149: swap Siria <-> Burkina Faso (name + values) ---
$ws.Range("A148").Value = "Burkina Faso"
$ws.Range("B148").Value = 1204
$ws.Range("C148").Value = 29
$ws.Range("D148").Value = 984
$ws.Range("E148").Value = 166
$ws.Range("H148").Value = 54

$ws.Range("A149").Value = "Siria"
$ws.Range("B149").Value = 1188
$ws.Range("C149").Value = 0
$ws.Range("D149").Value = 346
$ws.Range("E149").Value = 790
$ws.Range("H149").Value = 52

# --- Row 202/203: swap Santa Lucia <-> Timor Oriental (name only, values identical) ---
$ws.Range("A202").Value = "Timor Oriental"
$ws.Range("A203").Value = "Santa Lucia"

# --- Row 213/214: swap Islas Malvinas <-> Montserrat (name + values) ---
$ws.Range("A213").Value = "Montserrat"
$ws.Range("D213").Value = 12
$ws.Range("H213").Value = 1

$ws.Range("A214").Value = "Islas Malvinas"
$ws.Range("D214").Value = 13
$ws.Range("H214").Value = 0
